$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 ("Ministerios con servicios" gains a
# "Nivel Central" entry as its own row, shifting everything else down).
$ws.Rows("2:2").Insert()

# The inserted row copies formatting from the row above (the bold header
# row) - strip that back to the plain/no-style formatting used by all the
# other data rows.
$ws.Rows("2:2").ClearFormats()

# Fill every column (A:W) of the new row with "Nivel Central".
$ws.Range("A2:W2").Value = "Nivel Central"

# Restore the active-cell selection to where the author left off editing.
$ws.Range("B19").Select()
